$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 14212.5
$ws.Range("J17").Value = 16000
$ws.Range("L17").Value = 48000
$ws.Range("N17").Value = -48336
$ws.Range("H76").Value = 8862.875
$ws.Range("I76").Value = 8299.666999999999
$ws.Range("K76").Value = 8299.666999999999
$ws.Range("M76").Value = -7984.666999999999
$ws.Range("H79").Value = 8862.875
$ws.Range("I79").Value = 8299.666999999999
$ws.Range("K79").Value = 8299.666999999999
$ws.Range("M79").Value = -7207.666999999999
$ws.Range("H86").Value = 4016.7778
$ws.Range("J86").Value = 4091.8333
$ws.Range("L86").Value = 4091.8333
$ws.Range("N86").Value = -6337.8333
$ws.Range("H89").Value = 4016.7778
$ws.Range("J89").Value = 4091.8333
$ws.Range("L89").Value = 20459.1665
$ws.Range("N89").Value = -31691.1665
$ws.Range("H106").Value = 10783.777
$ws.Range("I106").Value = 4447
$ws.Range("K106").Value = 4447
$ws.Range("M106").Value = -3816
$ws.Range("H112").Value = 1791.7
$ws.Range("J112").Value = 1800.2106
$ws.Range("L112").Value = 5400.6318
$ws.Range("N112").Value = -7616.6318
$ws.Range("H128").Value = 36297.855
$ws.Range("J128").Value = 32986
$ws.Range("L128").Value = 32986
$ws.Range("N128").Value = -42946
$ws.Range("H132").Value = 1439.7241
$ws.Range("I132").Value = 1439.7241
$ws.Range("K132").Value = 4319.1723
$ws.Range("M132").Value = -1789.1723
$ws.Range("H137").Value = 13160754
$ws.Range("I137").Value = 52633924
$ws.Range("J137").Value = 3031.1228
$ws.Range("K137").Value = 157901772
$ws.Range("L137").Value = 9093.368399999999
$ws.Range("M137").Value = -157899222
$ws.Range("N137").Value = -14193.3684
$ws.Range("H138").Value = 3812.1936
$ws.Range("J138").Value = 4230.4614
$ws.Range("L138").Value = 12691.3842
$ws.Range("N138").Value = -22971.3842

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 741.3333
$ws.Range("I16").Value = 562
$ws.Range("K16").Value = 562
$ws.Range("M16").Value = -275
$ws.Range("H74").Value = 3274.55
$ws.Range("I74").Value = 2233.9412
$ws.Range("K74").Value = 2233.9412
$ws.Range("M74").Value = -1359.9412
$ws.Range("H77").Value = 3274.55
$ws.Range("I77").Value = 2233.9412
$ws.Range("K77").Value = 11169.706
$ws.Range("M77").Value = -6801.706000000002
$ws.Range("H97").Value = 1417.0714
$ws.Range("I97").Value = 1060.1111
$ws.Range("K97").Value = 1060.1111
$ws.Range("M97").Value = -564.1111000000001
$ws.Range("H132").Value = 3919.8064
$ws.Range("I132").Value = 2700.0454
$ws.Range("K132").Value = 8100.1362
$ws.Range("M132").Value = -5570.1362

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 21085.941
$ws.Range("I105").Value = 30698
$ws.Range("J105").Value = 14357.5
$ws.Range("K105").Value = 30698
$ws.Range("L105").Value = 14357.5
$ws.Range("M105").Value = -28951
$ws.Range("N105").Value = -17851.5
$ws.Range("H134").Value = 2315.6978
$ws.Range("I134").Value = 1374.0938
$ws.Range("J134").Value = 5054.909
$ws.Range("K134").Value = 4122.2814
$ws.Range("L134").Value = 15164.727
$ws.Range("M134").Value = -1587.2814
$ws.Range("N134").Value = -20234.727

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 170.75
$ws.Range("J38").Value = 116.666664
$ws.Range("L38").Value = 349.999992
$ws.Range("N38").Value = -1043.999992
$ws.Range("H113").Value = 142859620
$ws.Range("I113").Value = 2999.6667
$ws.Range("K113").Value = 8999.000100000001
$ws.Range("M113").Value = -6829.000100000001
$ws.Range("H124").Value = 41671016
$ws.Range("J124").Value = 66671330
$ws.Range("L124").Value = 200013990
$ws.Range("N124").Value = -200023810
$ws.Range("H137").Value = 2326.6667
$ws.Range("J137").Value = 4387.75
$ws.Range("L137").Value = 13163.25
$ws.Range("N137").Value = -23363.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""
$ws.Range("H62").Value = 44973.332
$ws.Range("I62").Value = 44973.332
$ws.Range("K62").Value = 44973.332
$ws.Range("M62").Value = -44287.332
$ws.Range("H65").Value = 44973.332
$ws.Range("I65").Value = 44973.332
$ws.Range("K65").Value = 134919.996
$ws.Range("M65").Value = -131487.996
$ws.Range("H80").Value = 253833.84
$ws.Range("I80").Value = 359307.5
$ws.Range("J80").Value = 7728.6665
$ws.Range("K80").Value = 359307.5
$ws.Range("L80").Value = 7728.6665
$ws.Range("M80").Value = -358309.5
$ws.Range("N80").Value = -9724.666499999999
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H83").Value = 253833.84
$ws.Range("I83").Value = 359307.5
$ws.Range("J83").Value = 7728.6665
$ws.Range("K83").Value = 1796537.5
$ws.Range("L83").Value = 38643.3325
$ws.Range("M83").Value = -1791545.5
$ws.Range("N83").Value = -48627.3325
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H113").Value = 3670.4707
$ws.Range("I113").Value = 3058.2856
$ws.Range("J113").Value = 4099
$ws.Range("K113").Value = 3058.2856
$ws.Range("L113").Value = 4099
$ws.Range("M113").Value = -888.2856000000002
$ws.Range("N113").Value = -8439
$ws.Range("H132").Value = 2255.2546
$ws.Range("I132").Value = 1791.22
$ws.Range("K132").Value = 5373.66
$ws.Range("M132").Value = -2843.66

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2427315.8
$ws.Range("I40").Value = 3251770.8
$ws.Range("J40").Value = 8914.532999999999
$ws.Range("K40").Value = 3251770.8
$ws.Range("L40").Value = 8914.532999999999
$ws.Range("M40").Value = -3251634.8
$ws.Range("N40").Value = -9186.532999999999
$ws.Range("H68").Value = 7636.0415
$ws.Range("I68").Value = 4133.5
$ws.Range("K68").Value = 4133.5
$ws.Range("M68").Value = -3384.5
$ws.Range("H71").Value = 7636.0415
$ws.Range("I71").Value = 4133.5
$ws.Range("K71").Value = 20667.5
$ws.Range("M71").Value = -16923.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2931
$ws.Range("I126").Value = 2318.3333
$ws.Range("J126").Value = 3599.3635
$ws.Range("K126").Value = 6954.999899999999
$ws.Range("L126").Value = 10798.0905
$ws.Range("M126").Value = -4484.999899999999
$ws.Range("N126").Value = -15738.0905
$ws.Range("H132").Value = 2629.8572
$ws.Range("I132").Value = 2125
$ws.Range("K132").Value = 6375
$ws.Range("M132").Value = -3845
